# adds ArtisanCommands moveBackground and pidLookahead
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row right after the "pidSource(<int>)" row (row 78) and
# fill it in with the new "pidLookahead" Artisan command.
$ws.Rows.Item(79).Insert()
$ws.Range("B79").Value = "pidLookahead(<int>)"
$ws.Range("C79").Value = "sets the PID lookahead"

# Insert a new row right after the "alarmset(<as>)" row (now row 90, since
# the previous insert shifted everything down by one) and fill it in with
# the new "moveBackground" Artisan command.
$ws.Rows.Item(91).Insert()
$ws.Range("B91").Value = "moveBackground(<direction>,<int>)"
$ws.Range("C91").Value = "moves the background profile the indicated number of steps towards <direction>, with <direction> one of up, down, left, right"

# Leave the selection on the newly added pidLookahead row, matching the
# state the workbook was left in.
$ws.Range("B79:C79").Select()
